$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" for the 54305882-... file
# (shared by Overview!G and de-de!H for the same source file rows)
$wsOverview.Range("G2").Value = "2016-09-07 06:22:54"
$wsOverview.Range("G5").Value = "2016-09-07 06:22:54"
$wsDeDe.Range("H2").Value = "2016-09-07 06:22:54"
$wsDeDe.Range("H5").Value = "2016-09-07 06:22:54"

# Priority changed from "ht" (human translation) to "mt" (machine translation)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# zh-cn Correspond Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-09-07 06:22:49"
$wsZhCn.Range("H5").Value = "2016-09-07 06:22:49"

# zh-cn Correspond Handback DateTime
$wsZhCn.Range("K2").Value = "2016-09-07 06:23:20"
$wsZhCn.Range("K5").Value = "2016-09-07 06:23:20"

# de-de Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-09-07 06:23:28"
$wsDeDe.Range("K5").Value = "2016-09-07 06:23:28"
